$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.086.83"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "3.059.76"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "386.69"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "102.37"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "36.63"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "0.0848"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "3.562.72"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").Value = "18.29"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "7.67"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "3.051.65"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").Value = "0.984"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "10.65"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "51.160.93"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  +3.21%  "
$ws.Range("D21").Value = "12.25"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "69.69"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "263.94"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "3.12"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").Value = "7.88"
$ws.Range("E26").Value = "  -4.73%  "
$ws.Range("D27").Value = "26.93"
$ws.Range("E27").Value = "  +3.09%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  -4.86%  "
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").Value = "10.39"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").Value = "35.36"
$ws.Range("E33").Value = "  +4.65%  "
$ws.Range("D34").Value = "0.0471"
$ws.Range("E34").Value = "  +5.22%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "49.96"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "3.37"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("D39").Value = "0.287"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "129.24"
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.115"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "16.39"
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("D44").Value = "3.80"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").Value = "21.61"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").Value = "2.47"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("D48").Value = "2.08"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "2.058.44"
$ws.Range("D50").Value = "9.47"
$ws.Range("E50").Value = "  +12.37%  "
$ws.Range("D51").Value = "0.924"
$ws.Range("E51").Value = "  +17.61%  "
